$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column K (rows 2-29) currently holds the observable name "sig_Z_eta".
# Rename it to "Z_rap" for every data row.
for ($r = 2; $r -le 29; $r++) {
    $ws.Cells.Item($r, 11).Value = "Z_rap"
}

# Update the active selection to match the saved view state (K31).
$ws.Range("K31").Select()
